# Corrección a Diebold Mariano y revisión de Cap1
# Updates the P_valores and Estadisticos_DM matrices with corrected values.

$wb = $excel.ActiveWorkbook

# --- P_valores sheet ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.05432601270742476
$wsP.Range("D2").Value = 0.08681744053941398
$wsP.Range("E2").Value = 0.1960191637946296
$wsP.Range("F2").Value = 0.1373715580499015

$wsP.Range("B3").Value = 0.05432601270742476
$wsP.Range("D3").Value = 0.6188064255278101
$wsP.Range("E3").Value = 0.4328186307396327
$wsP.Range("F3").Value = 0.8666685082759344

$wsP.Range("B4").Value = 0.08681744053941398
$wsP.Range("C4").Value = 0.6188064255278101
$wsP.Range("E4").Value = 0.2964603087371764
$wsP.Range("F4").Value = 0.5500793610103263

$wsP.Range("B5").Value = 0.1960191637946296
$wsP.Range("C5").Value = 0.4328186307396327
$wsP.Range("D5").Value = 0.2964603087371764
$wsP.Range("F5").Value = 0.5168672950651949

$wsP.Range("B6").Value = 0.1373715580499015
$wsP.Range("C6").Value = 0.8666685082759344
$wsP.Range("D6").Value = 0.5500793610103263
$wsP.Range("E6").Value = 0.5168672950651949

# --- Estadisticos_DM sheet ---
$wsE = $wb.Worksheets.Item("Estadisticos_DM")

$wsE.Range("C2").Value = -2.032696614884332
$wsE.Range("D2").Value = -1.792494279120127
$wsE.Range("E2").Value = -1.333475365443387
$wsE.Range("F2").Value = -1.541853417208297

$wsE.Range("B3").Value = 2.032696614884332
$wsE.Range("D3").Value = -0.5046727971642586
$wsE.Range("E3").Value = 0.7990242062687005
$wsE.Range("F3").Value = 0.1698651840618371

$wsE.Range("B4").Value = 1.792494279120127
$wsE.Range("C4").Value = 0.5046727971642586
$wsE.Range("E4").Value = 1.069451547175008
$wsE.Range("F4").Value = 0.6069776315826171

$wsE.Range("B5").Value = 1.333475365443387
$wsE.Range("C5").Value = -0.7990242062687005
$wsE.Range("D5").Value = -1.069451547175008
$wsE.Range("F5").Value = -0.6587894715330183

$wsE.Range("B6").Value = 1.541853417208297
$wsE.Range("C6").Value = -0.1698651840618371
$wsE.Range("D6").Value = -0.6069776315826171
$wsE.Range("E6").Value = 0.6587894715330183
